$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: section "3. Quản lý mượn trả sách" — team size " (2 người)" -> " (1 người)"
# ---------------------------------------------------------------------------
$rng1 = $d.Content
$found1 = $rng1.Find.Execute(" (2 người)", $true, $false, $false, $false, $false, $true, 1, $false, " (1 người)", 2)

# ---------------------------------------------------------------------------
# Change 2: section "6. Quản lý danh mục" — heading becomes bold, and the
# trailing " (Optional)" note becomes " (1/2 người)" (bold).
# ---------------------------------------------------------------------------
$rngHeading = $d.Content
$foundHeading = $rngHeading.Find.Execute("6. Quản lý danh mục", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($foundHeading) {
    $rngHeading.Bold = 1
    $headingEnd = $rngHeading.End

    # Narrow the search window to just after the heading so we only touch
    # the "(Optional)" that immediately follows it (the word also appears
    # after sections 7 and 8).
    $rngNote = $d.Range($headingEnd, $headingEnd + 60)
    $foundNote = $rngNote.Find.Execute(" (Optional)", $true, $false, $false, $false, $false, $true, 1, $false, " (1/2 người)", 2)
    if ($foundNote) {
        $rngNote.Bold = 1
    }
}

# ---------------------------------------------------------------------------
# Change 3: the "GUI" bullet — "(1/2 người)" -> "(1 người)"
# ---------------------------------------------------------------------------
$rngGui = $d.Content
$foundGui = $rngGui.Find.Execute("GUI ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($foundGui) {
    $guiEnd = $rngGui.End
    $rngGuiNote = $d.Range($guiEnd, $guiEnd + 40)
    $foundGuiNote = $rngGuiNote.Find.Execute("(1/2 người)", $true, $false, $false, $false, $false, $true, 1, $false, "(1 người)", 2)
}
